$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 29 de Julio de 2020 a las 18:14'

$ws.Range("B4").Value = 4517074
$ws.Range("C4").Value = 18731
$ws.Range("D4").Value = 2191005
$ws.Range("E4").Value = 2173323
$ws.Range("G4").Value = 405
$ws.Range("H4").Value = 152746

$ws.Range("B6").Value = 1579240
$ws.Range("C6").Value = 47105
$ws.Range("D6").Value = 1017204
$ws.Range("E6").Value = 527088
$ws.Range("G6").Value = 724
$ws.Range("H6").Value = 34948

$ws.Range("B11").Value = 351575
$ws.Range("C11").Value = 1775
$ws.Range("D11").Value = 324557
$ws.Range("E11").Value = 17740
$ws.Range("G11").Value = 38
$ws.Range("H11").Value = 9278

$ws.Range("B13").Value = 301455
$ws.Range("C13").Value = 763
$ws.Range("G13").Value = 83
$ws.Range("H13").Value = 45961

$ws.Range("B18").Value = 246776
$ws.Range("C18").Value = 289
$ws.Range("D18").Value = 199031
$ws.Range("E18").Value = 12616
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 35129

$ws.Range("B21").Value = 208537
$ws.Range("C21").Value = 586
$ws.Range("E21").Value = 7326
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 9211

$ws.Range("B25").Value = 115246
$ws.Range("C25").Value = 252
$ws.Range("D25").Value = 100308
$ws.Range("E25").Value = 6025
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 8913

$ws.Range("B33").Value = 79782
$ws.Range("C33").Value = 27
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 5730

$ws.Range("D45").Value = 46098
$ws.Range("E45").Value = 5406

$ws.Range("A72").Value = 'Chequia'
$ws.Range("B72").Value = 15986
$ws.Range("C72").Value = 434
$ws.Range("D72").Value = 11428
$ws.Range("E72").Value = 4184
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 374

$ws.Range("A73").Value = 'El Salvador'
$ws.Range("B73").Value = 15841
$ws.Range("C73").Value = 395
$ws.Range("D73").Value = 8071
$ws.Range("E73").Value = 7340
$ws.Range("G73").Value = 13
$ws.Range("H73").Value = 430

$ws.Range("A74").Value = 'Etiopia'
$ws.Range("B74").Value = 15810
$ws.Range("C74").Value = 610
$ws.Range("D74").Value = 6685
$ws.Range("E74").Value = 8872
$ws.Range("G74").Value = 14
$ws.Range("H74").Value = 253

$ws.Range("A75").Value = 'Costa de Marfil'
$ws.Range("B75").Value = 15713
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 10537
$ws.Range("E75").Value = 5078
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 98

$ws.Range("A76").Value = 'Australia'
$ws.Range("B76").Value = 15580
$ws.Range("C76").Value = 276
$ws.Range("D76").Value = 9431
$ws.Range("E76").Value = 5973
$ws.Range("G76").Value = 9
$ws.Range("H76").Value = 176

$ws.Range("A83").Value = 'Republica de Macedonia'
$ws.Range("B83").Value = 10503
$ws.Range("C83").Value = 188
$ws.Range("D83").Value = 5931
$ws.Range("E83").Value = 4096
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 476

$ws.Range("A84").Value = 'Madagascar'
$ws.Range("B84").Value = 10317
$ws.Range("C84").Value = 213
$ws.Range("D84").Value = 7117
$ws.Range("E84").Value = 3101
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 99

$ws.Range("A98").Value = 'Albania'
$ws.Range("B98").Value = 5105
$ws.Range("C98").Value = 108
$ws.Range("D98").Value = 2830
$ws.Range("E98").Value = 2125
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 150

$ws.Range("A99").Value = 'Republica de Yibuti'
$ws.Range("B99").Value = 5081
$ws.Range("C99").Value = 13
$ws.Range("D99").Value = 4999
$ws.Range("E99").Value = 24
$ws.Range("H99").Value = 58

$ws.Range("B116").Value = 2905
$ws.Range("C116").Value = 5
$ws.Range("E116").Value = 194
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 39

$ws.Range("B120").Value = 2521
$ws.Range("C120").Value = 1
$ws.Range("D120").Value = 1927
$ws.Range("E120").Value = 470

$ws.Range("B134").Value = 1748
$ws.Range("C134").Value = 28
$ws.Range("D134").Value = 616
$ws.Range("E134").Value = 1121

$ws.Range("A180").Value = 'Trinidad yTobago'
$ws.Range("C180").Value = 3
$ws.Range("D180").Value = 128
$ws.Range("E180").Value = 20
$ws.Range("H180").Value = 8

$ws.Range("A181").Value = 'Bermudas'
$ws.Range("B181").Value = 156
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 141
$ws.Range("E181").Value = 6
$ws.Range("H181").Value = 9
